$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 101654.555
$ws.Range("J17").Value = 101654.555
$ws.Range("L17").Value = 304963.665
$ws.Range("N17").Value = -305299.665
$ws.Range("H18").Value = 199333.33
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H28").Value = 265
$ws.Range("I28").Value = 271.72223
$ws.Range("K28").Value = 271.72223
$ws.Range("M28").Value = 213.27777
$ws.Range("H40").Value = 4214.143
$ws.Range("I40").Value = 3500
$ws.Range("J40").Value = 4499.8
$ws.Range("K40").Value = 3500
$ws.Range("L40").Value = 4499.8
$ws.Range("M40").Value = -3325
$ws.Range("N40").Value = -4849.8
$ws.Range("H53").Value = 342.34616
$ws.Range("I53").Value = 191.14285
$ws.Range("J53").Value = 518.75
$ws.Range("K53").Value = 191.14285
$ws.Range("L53").Value = 518.75
$ws.Range("M53").Value = 445.85715
$ws.Range("N53").Value = -1792.75
$ws.Range("H55").Value = 369.86667
$ws.Range("I55").Value = 149.5
$ws.Range("J55").Value = 516.7778
$ws.Range("K55").Value = 149.5
$ws.Range("L55").Value = 516.7778
$ws.Range("M55").Value = 64.5
$ws.Range("N55").Value = -944.7778
$ws.Range("H64").Value = 4071.2856
$ws.Range("I64").Value = 3749.8333
$ws.Range("K64").Value = 3749.8333
$ws.Range("M64").Value = -3501.8333
$ws.Range("H67").Value = 4071.2856
$ws.Range("I67").Value = 3749.8333
$ws.Range("K67").Value = 3749.8333
$ws.Range("M67").Value = -2891.8333
$ws.Range("H70").Value = 1712.5
$ws.Range("I70").Value = 1800
$ws.Range("J70").Value = 1625
$ws.Range("K70").Value = 5400
$ws.Range("L70").Value = 4875
$ws.Range("M70").Value = -5130
$ws.Range("N70").Value = -5415
$ws.Range("H73").Value = 1712.5
$ws.Range("I73").Value = 1800
$ws.Range("J73").Value = 1625
$ws.Range("K73").Value = 5400
$ws.Range("L73").Value = 4875
$ws.Range("M73").Value = -4464
$ws.Range("N73").Value = -6747
$ws.Range("H82").Value = 2111
$ws.Range("I82").Value = 2111
$ws.Range("K82").Value = 6333
$ws.Range("M82").Value = -5927
$ws.Range("H85").Value = 2111
$ws.Range("I85").Value = 2111
$ws.Range("K85").Value = 6333
$ws.Range("M85").Value = -4929
$ws.Range("H103").Value = 38462052
$ws.Range("I103").Value = 399.66666
$ws.Range("J103").Value = 50000550
$ws.Range("K103").Value = 1198.99998
$ws.Range("L103").Value = 150001650
$ws.Range("M103").Value = -612.9999800000001
$ws.Range("N103").Value = -150002822
$ws.Range("H112").Value = 4256.8423
$ws.Range("J112").Value = 4275.294
$ws.Range("L112").Value = 12825.882
$ws.Range("N112").Value = -15041.882
$ws.Range("H137").Value = 4270.1943
$ws.Range("I137").Value = 1507.125
$ws.Range("K137").Value = 4521.375
$ws.Range("M137").Value = -1971.375
$ws.Range("H138").Value = 2642.8286
$ws.Range("J138").Value = 4107.125
$ws.Range("L138").Value = 12321.375
$ws.Range("N138").Value = -22601.375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30369.875
$ws.Range("I32").Value = 16373.363
$ws.Range("J32").Value = 184331.5
$ws.Range("K32").Value = 16373.363
$ws.Range("L32").Value = 184331.5
$ws.Range("M32").Value = -16086.363
$ws.Range("N32").Value = -184905.5
$ws.Range("H45").Value = 483745.66
$ws.Range("I45").Value = 920945.5600000001
$ws.Range("K45").Value = 920945.5600000001
$ws.Range("M45").Value = -920568.5600000001
$ws.Range("H122").Value = 1756.8462
$ws.Range("I122").Value = 1623.6842
$ws.Range("K122").Value = 4871.0526
$ws.Range("M122").Value = -2421.0526

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1010.0833
$ws.Range("I107").Value = 1046.375
$ws.Range("J107").Value = 937.5
$ws.Range("K107").Value = 1046.375
$ws.Range("L107").Value = 937.5
$ws.Range("M107").Value = 873.625
$ws.Range("N107").Value = -4777.5
$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 812.25
$ws.Range("J10").Value = 999
$ws.Range("L10").Value = 999
$ws.Range("N10").Value = -1277
$ws.Range("H14").Value = 1941.5
$ws.Range("I14").Value = 1729.8
$ws.Range("J14").Value = 3000
$ws.Range("K14").Value = 1729.8
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = -1559.8
$ws.Range("N14").Value = -3340
$ws.Range("H58").Value = 1215.4166
$ws.Range("I58").Value = 1216.909
$ws.Range("K58").Value = 1216.909
$ws.Range("M58").Value = -1013.909
$ws.Range("H132").Value = 1821.925
$ws.Range("I132").Value = 1645.8918
$ws.Range("J132").Value = 3993
$ws.Range("K132").Value = 4937.6754
$ws.Range("L132").Value = 11979
$ws.Range("M132").Value = -2407.6754
$ws.Range("N132").Value = -17039
$ws.Range("H136").Value = 1215.4166
$ws.Range("I136").Value = 1216.909
$ws.Range("K136").Value = 3650.727
$ws.Range("M136").Value = -1100.727

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 74202.21000000001
$ws.Range("I39").Value = 167496.83
$ws.Range("J39").Value = 4231.25
$ws.Range("K39").Value = 502490.49
$ws.Range("L39").Value = 12693.75
$ws.Range("M39").Value = -502196.49
$ws.Range("N39").Value = -13281.75
$ws.Range("H132").Value = 1289.8235
$ws.Range("I132").Value = 1103.0834
$ws.Range("J132").Value = 1738
$ws.Range("K132").Value = 9927.750599999999
$ws.Range("L132").Value = 15642
$ws.Range("M132").Value = -7397.750599999999
$ws.Range("N132").Value = -20702
$ws.Range("H138").Value = 4016.5
$ws.Range("J138").Value = 4016.5
$ws.Range("L138").Value = 12049.5
$ws.Range("N138").Value = -22329.5
$ws.Range("H141").Value = 2886.5454
$ws.Range("I141").Value = 2886.5454
$ws.Range("K141").Value = 8659.636200000001
$ws.Range("M141").Value = -3479.636200000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6221.8
$ws.Range("I70").Value = 5369.6665
$ws.Range("K70").Value = 5369.6665
$ws.Range("M70").Value = -5099.6665
$ws.Range("H73").Value = 6221.8
$ws.Range("I73").Value = 5369.6665
$ws.Range("K73").Value = 5369.6665
$ws.Range("M73").Value = -4433.6665
$ws.Range("H97").Value = 35357.74
$ws.Range("I97").Value = 36988.723
$ws.Range("K97").Value = 36988.723
$ws.Range("M97").Value = -36492.723

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19345.934
$ws.Range("I7").Value = 27118.9
$ws.Range("J7").Value = 3800
$ws.Range("K7").Value = 27118.9
$ws.Range("L7").Value = 3800
$ws.Range("M7").Value = -27006.9
$ws.Range("N7").Value = -4024
$ws.Range("H61").Value = 1893
$ws.Range("I61").Value = 1894.4615
$ws.Range("J61").Value = 1883.5
$ws.Range("K61").Value = 1894.4615
$ws.Range("L61").Value = 1883.5
$ws.Range("M61").Value = -1692.4615
$ws.Range("N61").Value = -2287.5
$ws.Range("H82").Value = 2259.8
$ws.Range("J82").Value = 2999.6667
$ws.Range("L82").Value = 2999.6667
$ws.Range("N82").Value = -3721.6667
$ws.Range("H85").Value = 2259.8
$ws.Range("J85").Value = 2999.6667
$ws.Range("L85").Value = 2999.6667
$ws.Range("N85").Value = -5495.6667
$ws.Range("H93").Value = 2262.2778
$ws.Range("I93").Value = 1921.9166
$ws.Range("J93").Value = 2943
$ws.Range("K93").Value = 1921.9166
$ws.Range("L93").Value = 2943
$ws.Range("M93").Value = -673.9166
$ws.Range("N93").Value = -5439
$ws.Range("H113").Value = 1893
$ws.Range("I113").Value = 1894.4615
$ws.Range("J113").Value = 1883.5
$ws.Range("K113").Value = 1894.4615
$ws.Range("L113").Value = 1883.5
$ws.Range("M113").Value = 275.5385000000001
$ws.Range("N113").Value = -6223.5
$ws.Range("H126").Value = 19345.934
$ws.Range("I126").Value = 27118.9
$ws.Range("J126").Value = 3800
$ws.Range("K126").Value = 81356.70000000001
$ws.Range("L126").Value = 11400
$ws.Range("M126").Value = -78886.70000000001
$ws.Range("N126").Value = -16340

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 40495
$ws.Range("J51").Value = 40495
$ws.Range("L51").Value = 40495
$ws.Range("N51").Value = -41515
$ws.Range("H96").Value = 4072.8333
$ws.Range("I96").Value = 4287.6
$ws.Range("J96").Value = 2999
$ws.Range("K96").Value = 4287.6
$ws.Range("L96").Value = 2999
$ws.Range("M96").Value = -2914.6
$ws.Range("N96").Value = -5745
$ws.Range("H107").Value = 83334700
$ws.Range("I107").Value = 1625
$ws.Range("J107").Value = 250000850
$ws.Range("K107").Value = 4875
$ws.Range("L107").Value = 750002550
$ws.Range("M107").Value = -2955
$ws.Range("N107").Value = -750006390
$ws.Range("H122").Value = 4276.522
$ws.Range("I122").Value = 3969.524
$ws.Range("K122").Value = 11908.572
$ws.Range("M122").Value = -9458.572
$ws.Range("H126").Value = 3295.0833
$ws.Range("J126").Value = 5124.75
$ws.Range("L126").Value = 15374.25
$ws.Range("N126").Value = -20314.25
$ws.Range("H136").Value = 3472.742
$ws.Range("I136").Value = 3846.4348
$ws.Range("K136").Value = 11539.3044
$ws.Range("M136").Value = -8989.304400000001
